$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current order (rows 16-18):
#  Row16: 1047364359 / LAURA RAQUEL ORTEGA ZAPATA      / F=5333  / G=1000000
#  Row17: 1007977741 / DIEGO ALEJANDRO NARVAEZ GARAY    / F=15733 / G=1475000
#  Row18: 22789716   / PATRICIA ISABEL TAMARA VARGAS    / F=5333  / G=1000000
#
# Target order (rows 16-18):
#  Row16: 22789716   / PATRICIA ISABEL TAMARA VARGAS    / F=5333  / G=1000000
#  Row17: 1047364359 / LAURA RAQUEL ORTEGA ZAPATA       / F=5333  / G=1000000
#  Row18: 1007977741 / DIEGO ALEJANDRO NARVAEZ GARAY    / F=15733 / G=2633409

$ws.Range("C16").Value = "22789716"
$ws.Range("D16").Value = "PATRICIA ISABEL TAMARA VARGAS"
$ws.Range("F16").Value = 5333
$ws.Range("G16").Value = 1000000

$ws.Range("C17").Value = "1047364359"
$ws.Range("D17").Value = "LAURA RAQUEL ORTEGA ZAPATA"
$ws.Range("F17").Value = 5333
$ws.Range("G17").Value = 1000000

$ws.Range("C18").Value = "1007977741"
$ws.Range("D18").Value = "DIEGO ALEJANDRO NARVAEZ GARAY"
$ws.Range("F18").Value = 15733
$ws.Range("G18").Value = 2633409

$wb.Save()
